# "Hjemme passive updated meanEMG legmaxROM"
# Update the AT_ind_max (legmaxROM/meanEMG) figures for subject columns B:E
# (rows 1-3: sample-size row, CON row, STR row) and re-select the range that
# actually holds the new data (B1:E3) instead of the old B1:AY3 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 - sample counts per column
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - CON
$ws.Range("B2").Value = 67.344776313292996
$ws.Range("C2").Value = 36.669401897941171
$ws.Range("D2").Value = 67.079255776434451
$ws.Range("E2").Value = 40.905546271713845

# Row 3 - STR
$ws.Range("B3").Value = 63.376823674849284
$ws.Range("C3").Value = 43.216688876332171
$ws.Range("D3").Value = 52.717170761474343
$ws.Range("E3").Value = 45.18835969066658

# Selection now reflects the updated (smaller) data range
$null = $ws.Range("B1:E3").Select()
